# Scheduled runner update: refresh cached market-board price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) on the per-job profit
# sheets. Values come from an external price feed snapshot; only the
# numeric market-data cells change, nothing structural.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 401028.8
$ws.Range("I28").Value = 500730.75
$ws.Range("K28").Value = 500730.75
$ws.Range("M28").Value = -500245.75

$ws.Range("H80").Value = 90919430
$ws.Range("I80").Value = 250000770
$ws.Range("J80").Value = 15807.286
$ws.Range("K80").Value = 750002310
$ws.Range("L80").Value = 47421.858
$ws.Range("M80").Value = -750001312
$ws.Range("N80").Value = -49417.858

$ws.Range("H83").Value = 90919430
$ws.Range("I83").Value = 250000770
$ws.Range("J83").Value = 15807.286
$ws.Range("K83").Value = 2250006930
$ws.Range("L83").Value = 142265.574
$ws.Range("M83").Value = -2250001938
$ws.Range("N83").Value = -152249.574

$ws.Range("H137").Value = 4701412
$ws.Range("I137").Value = 275793.34
$ws.Range("J137").Value = 6176618.5
$ws.Range("K137").Value = 827380.02
$ws.Range("L137").Value = 18529855.5
$ws.Range("M137").Value = -824830.02
$ws.Range("N137").Value = -18534955.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2814.6924
$ws.Range("I63").Value = 2026.7778
$ws.Range("K63").Value = 2026.7778
$ws.Range("M63").Value = -1340.7778

$ws.Range("H66").Value = 2814.6924
$ws.Range("I66").Value = 2026.7778
$ws.Range("K66").Value = 10133.889
$ws.Range("M66").Value = -6701.889000000001

$ws.Range("H74").Value = 2168.7585
$ws.Range("I74").Value = 2250.1
$ws.Range("K74").Value = 2250.1
$ws.Range("M74").Value = -1376.1

$ws.Range("H77").Value = 2168.7585
$ws.Range("I77").Value = 2250.1
$ws.Range("K77").Value = 11250.5
$ws.Range("M77").Value = -6882.5

$ws.Range("H97").Value = 1533.3334
$ws.Range("I97").Value = 1373.4
$ws.Range("K97").Value = 1373.4
$ws.Range("M97").Value = -877.4000000000001

$ws.Range("H113").Value = 113329.664
$ws.Range("J113").Value = 113329.664
$ws.Range("L113").Value = 113329.664
$ws.Range("N113").Value = -122007.664

$ws.Range("H122").Value = 3491.5173
$ws.Range("I122").Value = 2696.7144
$ws.Range("K122").Value = 8090.1432
$ws.Range("M122").Value = -5640.1432

$ws.Range("H132").Value = 403342.1
$ws.Range("I132").Value = 438067.75
$ws.Range("K132").Value = 1314203.25
$ws.Range("M132").Value = -1311673.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 29274.264
$ws.Range("I20").Value = 52691.8
$ws.Range("J20").Value = 3254.7778
$ws.Range("K20").Value = 52691.8
$ws.Range("L20").Value = 3254.7778
$ws.Range("M20").Value = -52444.8
$ws.Range("N20").Value = -3748.7778

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4364.58
$ws.Range("I31").Value = 1922.7222
$ws.Range("K31").Value = 1922.7222
$ws.Range("M31").Value = -1627.7222

$ws.Range("H34").Value = 4364.58
$ws.Range("I34").Value = 1922.7222
$ws.Range("K34").Value = 1922.7222
$ws.Range("M34").Value = -1720.7222

$ws.Range("H62").Value = 2631.9
$ws.Range("I62").Value = 2065.375
$ws.Range("J62").Value = 4898
$ws.Range("K62").Value = 2065.375
$ws.Range("L62").Value = 4898
$ws.Range("M62").Value = -1441.375
$ws.Range("N62").Value = -6146

$ws.Range("H65").Value = 2631.9
$ws.Range("I65").Value = 2065.375
$ws.Range("J65").Value = 4898
$ws.Range("K65").Value = 10326.875
$ws.Range("L65").Value = 24490
$ws.Range("M65").Value = -7206.875
$ws.Range("N65").Value = -30730

$ws.Range("H122").Value = 3529.68
$ws.Range("I122").Value = 2945
$ws.Range("J122").Value = 4163.0835
$ws.Range("K122").Value = 8835
$ws.Range("L122").Value = 12489.2505
$ws.Range("M122").Value = -6385
$ws.Range("N122").Value = -17389.2505

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 3000
$ws.Range("I70").Value = 3000
$ws.Range("K70").Value = 9000
$ws.Range("M70").Value = -8685

$ws.Range("H73").Value = 3000
$ws.Range("I73").Value = 3000
$ws.Range("K73").Value = 9000
$ws.Range("M73").Value = -7908

$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

$ws.Range("H103").Value = 371.85715
$ws.Range("I103").Value = 441.66666
$ws.Range("J103").Value = 319.5
$ws.Range("K103").Value = 1324.99998
$ws.Range("L103").Value = 958.5
$ws.Range("M103").Value = -445.9999800000001
$ws.Range("N103").Value = -2716.5

$ws.Range("H136").Value = 1455.7
$ws.Range("I136").Value = 1455.7
$ws.Range("K136").Value = 4367.1
$ws.Range("M136").Value = 732.8999999999996

$ws.Range("H138").Value = 21301238
$ws.Range("I138").Value = 1547.5
$ws.Range("K138").Value = 4642.5
$ws.Range("M138").Value = 497.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2985.111
$ws.Range("I80").Value = 2875.8333
$ws.Range("K80").Value = 2875.8333
$ws.Range("M80").Value = -1877.8333

$ws.Range("H83").Value = 2985.111
$ws.Range("I83").Value = 2875.8333
$ws.Range("K83").Value = 14379.1665
$ws.Range("M83").Value = -9387.166499999999

$ws.Range("H132").Value = 3843.4285
$ws.Range("I132").Value = 3501.4
$ws.Range("K132").Value = 10504.2
$ws.Range("M132").Value = -7974.200000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 2152.8948
$ws.Range("I55").Value = 2020.091
$ws.Range("J55").Value = 2335.5
$ws.Range("K55").Value = 2020.091
$ws.Range("L55").Value = 2335.5
$ws.Range("M55").Value = -1847.091
$ws.Range("N55").Value = -2681.5

$ws.Range("H68").Value = 1529.2
$ws.Range("I68").Value = 1544.8889
$ws.Range("J68").Value = 1505.6666
$ws.Range("K68").Value = 1544.8889
$ws.Range("L68").Value = 1505.6666
$ws.Range("M68").Value = -795.8888999999999
$ws.Range("N68").Value = -3003.6666

$ws.Range("H71").Value = 1529.2
$ws.Range("I71").Value = 1544.8889
$ws.Range("J71").Value = 1505.6666
$ws.Range("K71").Value = 7724.4445
$ws.Range("L71").Value = 7528.333000000001
$ws.Range("M71").Value = -3980.4445
$ws.Range("N71").Value = -15016.333

$ws.Range("H132").Value = 1001990.1
$ws.Range("I132").Value = 1001990.1
$ws.Range("K132").Value = 3005970.3
$ws.Range("M132").Value = -3003440.3

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 116662.664
$ws.Range("J16").Value = 116662.664
$ws.Range("L16").Value = 116662.664
$ws.Range("N16").Value = -117246.664

$ws.Range("H81").Value = 70815
$ws.Range("I81").Value = 145605.14
$ws.Range("J81").Value = 5373.625
$ws.Range("K81").Value = 291210.28
$ws.Range("L81").Value = 10747.25
$ws.Range("M81").Value = -290149.28
$ws.Range("N81").Value = -12869.25

$ws.Range("H84").Value = 70815
$ws.Range("I84").Value = 145605.14
$ws.Range("J84").Value = 5373.625
$ws.Range("K84").Value = 1456051.4
$ws.Range("L84").Value = 53736.25
$ws.Range("M84").Value = -1450747.4
$ws.Range("N84").Value = -64344.25

$ws.Range("H132").Value = 51291.285
$ws.Range("I132").Value = 61901.293
$ws.Range("J132").Value = 6198.75
$ws.Range("K132").Value = 185703.879
$ws.Range("L132").Value = 18596.25
$ws.Range("M132").Value = -183173.879
$ws.Range("N132").Value = -23656.25

$ws.Range("H136").Value = 3137.6924
$ws.Range("I136").Value = 1878.5
$ws.Range("K136").Value = 5635.5
$ws.Range("M136").Value = -3085.5
